# Updates the "NEW" sheet data table: two obsolete case rows were removed
# from the source export (Caso 6072 / PAZ SOLDAN 4991 and Caso -495 /
# Ricardo Balbin 3827), and every following row shifted up to close the gap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row that now holds Caso -495 (Ricardo Balbin 3827) first so the
# row-35 deletion below doesn't change its row number.
$ws.Rows("49:49").Delete()

# Remove the row that holds Caso 6072 (PAZ SOLDAN 4991); everything below
# shifts up by one, matching the target layout (data now ends at row 49).
$ws.Rows("35:35").Delete()
